$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2 through 15
# from serial date 45212 (2023-10-13) to serial date 45221 (2023-10-22).
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
